$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price values in column D (Price).
# Values are written with a leading apostrophe so Excel keeps them
# as text (matching the source data's string formatting, including
# trailing zeros), then the style is reset to Normal on each cell so
# no stray number-format is left applied.
$ws.Range("D2").Value = "'267.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.253"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06155"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.564"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.366"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8218"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.1558"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08191"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03304"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.03174"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09286"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.739"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001622"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04684"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006327"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005773"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001067"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.718"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.268"
$ws.Range("D24").Style = "Normal"
$ws.Range("D40").Value = "'0.04662"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007005"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.003895"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Value = "'0.01184"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006020"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.0009885"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Value = "'0.002438"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00001898"
$ws.Range("D50").Style = "Normal"
